{"js": "// Add the \"(max. N Zeichen)\" character-limit hints to the submission\n// template's prompt paragraphs (title + the various question prompts).\n\nconst body = context.document.body;\n\n// 1) \"Titel: Deep Learning ...\" -> \"Titel (max. 140 Zeichen): Deep Learning ...\"\n{\n  const found = body.search(\"Titel\", { matchCase: true, matchWholeWord: true });\n  found.load(\"items/text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\" (max. 140 Zeichen)\", \"After\");\n    await context.sync();\n  }\n}\n\n// 2) \"Was ist das Ziel des Projekts?\" -> add \" (max. 350 Zeichen)\"\n{\n  const found = body.search(\"Was ist das Ziel des Projekts?\", { matchCase: true });\n  found.load(\"items/text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\" (max. 350 Zeichen)\", \"After\");\n    await context.sync();\n  }\n}\n\n// 3) \"Warum ist euer Projekt wichtig?\" -> add \" (max. 280 Zeichen)\"\n{\n  const found = body.search(\"Warum ist euer Projekt wichtig?\", { matchCase: true });\n  found.load(\"items/text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\" (max. 280 Zeichen)\", \"After\");\n    await context.sync();\n  }\n}\n\n// 4) \"Wer kann eure Ergebnisse verwenden? Was ist der Anwendungsfall?\" -> add \" (max. 210 Zeichen)\"\n{\n  const found = body.search(\"Wer kann eure Ergebnisse verwenden? Was ist der Anwendungsfall?\", { matchCase: true });\n  found.load(\"items/text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\" (max. 210 Zeichen)\", \"After\");\n    await context.sync();\n  }\n}\n\n// 5) \"Beschreibung des Datensatzes:\" -> \"Beschreibung des Datensatzes (max. 350 Zeichen): \"\n{\n  const found = body.search(\"Beschreibung des Datensatzes:\", { matchCase: true });\n  found.load(\"items/text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\"Beschreibung des Datensatzes (max. 350 Zeichen): \", \"Replace\");\n    await context.sync();\n  }\n}\n\n// 6) \"Aufbereitung der Daten und sonstige Vorbereitung:\" -> \"... (max. 350 Zeichen):\"\n{\n  const found = body.search(\"Aufbereitung der Daten und sonstige Vorbereitung:\", { matchCase: true });\n  found.load(\"items/text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\"Aufbereitung der Daten und sonstige Vorbereitung (max. 350 Zeichen):\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// 7) \"...Trainieren und Testen der Daten...):\" -> \"...Daten...) (max. 350 Zeichen):\"\n{\n  const found = body.search(\"Trainieren und Testen der Daten...):\", { matchCase: true });\n  found.load(\"items/text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\"Trainieren und Testen der Daten...) (max. 350 Zeichen):\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Add the \"(max. N Zeichen)\" character-limit hints to the submission\n# template's prompt paragraphs (title + the various question prompts).\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n    $range = $d.Content\n    $range.Find.Execute(\n        $FindText,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $ReplaceText, # ReplaceWith\n        2            # Replace (wdReplaceOne)\n    ) | Out-Null\n}\n\n# 1) \"Titel: Deep Learning ...\" -> \"Titel (max. 140 Zeichen): Deep Learning ...\"\nReplace-DocText \"Titel: Deep Learning zur Grundrissplan-Analyse - Detektion von Objektsymbolen\" \"Titel (max. 140 Zeichen): Deep Learning zur Grundrissplan-Analyse - Detektion von Objektsymbolen\"\n\n# 2) \"Was ist das Ziel des Projekts?\" -> add \" (max. 350 Zeichen)\"\nReplace-DocText \"Was ist das Ziel des Projekts?\" \"Was ist das Ziel des Projekts? (max. 350 Zeichen)\"\n\n# 3) \"Warum ist euer Projekt wichtig?\" -> add \" (max. 280 Zeichen)\"\nReplace-DocText \"Warum ist euer Projekt wichtig?\" \"Warum ist euer Projekt wichtig? (max. 280 Zeichen)\"\n\n# 4) \"Wer kann eure Ergebnisse verwenden? Was ist der Anwendungsfall?\" -> add \" (max. 210 Zeichen)\"\nReplace-DocText \"Wer kann eure Ergebnisse verwenden? Was ist der Anwendungsfall?\" \"Wer kann eure Ergebnisse verwenden? Was ist der Anwendungsfall? (max. 210 Zeichen)\"\n\n# 5) \"Beschreibung des Datensatzes:\" -> \"Beschreibung des Datensatzes (max. 350 Zeichen): \"\nReplace-DocText \"Beschreibung des Datensatzes:\" \"Beschreibung des Datensatzes (max. 350 Zeichen): \"\n\n# 6) \"Aufbereitung der Daten und sonstige Vorbereitung:\" -> \"... (max. 350 Zeichen):\"\nReplace-DocText \"Aufbereitung der Daten und sonstige Vorbereitung:\" \"Aufbereitung der Daten und sonstige Vorbereitung (max. 350 Zeichen):\"\n\n# 7) \"...Trainieren und Testen der Daten...):\" -> \"...Daten...) (max. 350 Zeichen):\"\nReplace-DocText \"Trainieren und Testen der Daten...):\" \"Trainieren und Testen der Daten...) (max. 350 Zeichen):\"\n"}
